$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.753.99"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "1.642.59"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.42"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("E6").Value = "  -0.57%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0630"
$ws.Range("E8").Value = "  +0.62%  "
$ws.Range("E9").Value = "  -0.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.15"
$ws.Range("E10").Value = "  -0.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0842"
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("D12").Value = "1.866.89"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").Value = "1.647.49"
$ws.Range("E13").Value = "  +2.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.16"
$ws.Range("E14").Value = "  -1.33%  "
$ws.Range("E15").Value = "  -1.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.35"
$ws.Range("E16").Value = "  -2.46%  "
$ws.Range("D17").Value = "26.744.26"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("E18").Value = "  -1.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "213.94"
$ws.Range("E19").Value = "  -1.84%  "
$ws.Range("E20").Value = "  +0.29%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("E22").Value = "  +13.91%  "
$ws.Range("E23").Value = "  -0.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.34"
$ws.Range("E24").Value = "  -2.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.89"
$ws.Range("E25").Value = "  -1.18%  "
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("E27").Value = "  -1.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.10"
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.63"
$ws.Range("E29").Value = "  -1.41%  "
$ws.Range("E30").Value = "  -1.37%  "
$ws.Range("E31").Value = "  +0.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.32"
$ws.Range("E32").Value = "  -1.99%  "
$ws.Range("E33").Value = "  -2.06%  "
$ws.Range("D34").Value = "1.293.91"
$ws.Range("E34").Value = "  +1.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.54"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.44"
$ws.Range("E36").Value = "  +1.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0174"
$ws.Range("E37").Value = "  -4.88%  "
$ws.Range("E38").Value = "  +1.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.826"
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("E40").Value = "  +0.25%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").Value = "  -0.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.35"
$ws.Range("E43").Value = "  -1.90%  "
$ws.Range("D44").Value = "1.792.84"
$ws.Range("E44").Value = "  +0.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.31"
$ws.Range("E45").Value = "  +2.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.43"
$ws.Range("E46").Value = "  -1.86%  "
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0101"
$ws.Range("E49").Value = "  -4.71%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.70"
$ws.Range("E50").Value = "  -1.29%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0976"
$ws.Range("E51").Value = "  -0.03%  "
